$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("好聽的" / "清唱" move from C to D)
$ws.Range("C1").EntireColumn.Insert()

# Fill the newly inserted column C with what used to be in column B ("日文")
$ws.Range("C1").Value = "日文"
$ws.Range("C2").Value = "日文"

# Column B now holds the new "lian" values
$ws.Range("B1").Value = "lian"
$ws.Range("B2").Value = "lian"

# Move the active selection to C5, matching the edited workbook's view state
$ws.Range("C5").Select()
